# update with new logo and colors
# (content/metadata refresh of the Metadata sheet: version bump, status
# change, new date, updated contact info, new Jurisdiction row)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- simple in-place field updates -----------------------------------
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- make room for a new "Jurisdiction" row at row 12 -----------------
# Push the existing Description / Purpose / Copyright / Immutable block
# down by one row (bottom-up so we never clobber a value before reading
# it), carrying the row formatting along with it.

$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A13:B13").Copy()
$ws.Range("A14:B14").PasteSpecial(-4122)
$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A12:B12").Copy()
$ws.Range("A13:B13").PasteSpecial(-4122)
$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "RxNorm codes for all of Form 2400"

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

$excel.CutCopyMode = $false
